# Rename the sole worksheet from "Roma DST Perm and Pressure 2012" to "data".
# Renaming via the Name property also updates any formulas / defined names
# (e.g. the hidden _xlnm._FilterDatabase name) that reference the sheet by
# name, matching the workbook.xml diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "data"

# Re-assert the existing selection (B14) so it is left untouched, then scroll
# the window so row 189 is at the top of the visible pane (sheetView
# topLeftCell="A189" in the diff) without disturbing the active cell.
$ws.Range("B14").Select()
$excel.ActiveWindow.ScrollRow = 189
